# 20210917 update codes from Google Drive (lastest 20210917)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: title / teacher / semester / year
$ws.Range("A1").Value = "ทดสอบ"
$ws.Range("B1").Value = "คุณครู"
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2021

# Row 2: column headers (added "ชื่อระดับการศึกษา" before "ระดับชั้นเรียน")
$ws.Range("A2").Value = "date"
$ws.Range("B2").Value = "เวลาในการสอน"
$ws.Range("C2").Value = "ชั้น"
$ws.Range("D2").Value = "รหัสวิชา"
$ws.Range("E2").Value = "ชื่อระดับการศึกษา"
$ws.Range("F2").Value = "ระดับชั้นเรียน"

# Row 3: data values updated to new class/schedule
$ws.Range("A3").Value = "วันศุกร์"
$ws.Range("B3").Value = "15:00 - 16:00"
$ws.Range("C3").Value = "ป.4/2"
$ws.Range("D3").Value = "ส 21101"
$ws.Range("E3").Value = "ประถมศึกษา"
$ws.Range("F3").Value = "ประถมปลาย"

# Column widths (closest achievable values to the target OOXML widths)
$ws.Columns.Item(1).ColumnWidth = 7.333333
$ws.Columns.Item(2).ColumnWidth = 13.666667
$ws.Columns.Item(3).ColumnWidth = 9.833333
$ws.Columns.Item(4).ColumnWidth = 14
$ws.Columns.Item(5).ColumnWidth = 15.666667
$ws.Columns.Item(6).ColumnWidth = 11.333333

# Selection moved to E15
$ws.Range("E15").Select() | Out-Null
